$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Sheet1: Forecast Comparison ---
$ws1.Range("D2").Value = 175
$ws1.Range("H2").Value = 7.35
$ws1.Range("L2").Value = 1.02
$ws1.Range("D3").Value = 177
$ws1.Range("H3").Value = 6.3
$ws1.Range("L3").Value = 0.88
$ws1.Range("D4").Value = 175
$ws1.Range("H4").Value = 5.34
$ws1.Range("L4").Value = 1.12
$ws1.Range("D5").Value = 173
$ws1.Range("H5").Value = 4.39
$ws1.Range("L5").Value = 1
$ws1.Range("D6").Value = 174
$ws1.Range("H6").Value = 3.37
$ws1.Range("L6").Value = 1.1
$ws1.Range("D7").Value = 168
$ws1.Range("H7").Value = 2.45
$ws1.Range("L7").Value = 0.88
$ws1.Range("D8").Value = 175
$ws1.Range("H8").Value = 1.4
$ws1.Range("L8").Value = 1.12
$ws1.Range("D9").Value = 148
$ws1.Range("H9").Value = 0.47
$ws1.Range("I9").Value = "High"
$ws1.Range("L9").Value = 1.11
$ws1.Range("D10").Value = 145
$ws1.Range("L10").Value = 1.04
$ws1.Range("D11").Value = 146
$ws1.Range("L11").Value = 1.03
$ws1.Range("D12").Value = 174
$ws1.Range("L12").Value = 1.17
$ws1.Range("L13").Value = 1.09
$ws1.Range("L14").Value = 1.15
$ws1.Range("D15").Value = 142
$ws1.Range("L15").Value = 1.12
$ws1.Range("D16").Value = 143
$ws1.Range("L16").Value = 1.11
$ws1.Range("D17").Value = 143
$ws1.Range("L17").Value = 0.9

# --- Sheet2: Summary (values stored as text) ---
$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "2608"
$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "1369"
$ws2.Range("B11").NumberFormat = "@"
$ws2.Range("B11").Value = "702"
$ws2.Range("B12").NumberFormat = "@"
$ws2.Range("B12").Value = "177"
$ws2.Range("B14").NumberFormat = "@"
$ws2.Range("B14").Value = "143"
